$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds numeric-looking price strings that must stay text (t="inlineStr" in
# the source). Force text via NumberFormat="@" while assigning, then restore the
# cell style to Normal so no extra formatting is left applied to the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "250.54"
Set-TextValue $ws.Cells.Item(3, 4) "23.37"
$ws.Cells.Item(6, 2).Value = "GateToken"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Cells.Item(6, 4) "3.427"
$ws.Cells.Item(6, 5).Value = "5GateTokenGT"
$ws.Cells.Item(7, 2).Value = "KuCoinToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Cells.Item(7, 4) "6.570"
$ws.Cells.Item(7, 5).Value = "6KuCoinTokenKCS"
Set-TextValue $ws.Cells.Item(8, 4) "1.330"
Set-TextValue $ws.Cells.Item(9, 4) "0.7930"
$ws.Cells.Item(10, 2).Value = "WazirX"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Cells.Item(10, 4) "0.1486"
$ws.Cells.Item(10, 5).Value = "9WazirXWRX"
$ws.Cells.Item(11, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Cells.Item(11, 4) "0.07830"
$ws.Cells.Item(11, 5).Value = "10MandalaExchangeTokenMDX"
$ws.Cells.Item(12, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Cells.Item(12, 4) "0.03348"
$ws.Cells.Item(12, 5).Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Cells.Item(13, 4) "0.03035"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"
$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Cells.Item(14, 4) "0.09257"
$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"
$ws.Cells.Item(15, 2).Value = "MCDex"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Cells.Item(15, 4) "3.566"
$ws.Cells.Item(15, 5).Value = "14MCDexMCB"
$ws.Cells.Item(16, 2).Value = "BitForexToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Cells.Item(16, 4) "0.001687"
$ws.Cells.Item(16, 5).Value = "15BitForexTokenBF"
$ws.Cells.Item(17, 2).Value = "CoinExToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Cells.Item(17, 4) "0.04783"
$ws.Cells.Item(17, 5).Value = "16CoinExTokenCET"
$ws.Cells.Item(18, 2).Value = "One"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Cells.Item(18, 4) "0.0006063"
$ws.Cells.Item(18, 5).Value = "17OneONE"
Set-TextValue $ws.Cells.Item(19, 4) "0.006201"
Set-TextValue $ws.Cells.Item(20, 4) "0.005682"
Set-TextValue $ws.Cells.Item(21, 4) "0.001068"
Set-TextValue $ws.Cells.Item(22, 4) "0.0001501"
Set-TextValue $ws.Cells.Item(23, 4) "3.695"
Set-TextValue $ws.Cells.Item(26, 4) "0.1256"
Set-TextValue $ws.Cells.Item(27, 4) "0.0006479"
Set-TextValue $ws.Cells.Item(40, 4) "0.04438"
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Cells.Item(42, 4) "0.1067"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Cells.Item(43, 4) "0.003152"
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"
Set-TextValue $ws.Cells.Item(44, 4) "0.009246"
$ws.Cells.Item(45, 5).Value = "44ACDXExchangeACXTBestin24h"
Set-TextValue $ws.Cells.Item(46, 4) "0.00005897"
Set-TextValue $ws.Cells.Item(48, 4) "0.9907"
Set-TextValue $ws.Cells.Item(49, 4) "0.1040"
Set-TextValue $ws.Cells.Item(51, 4) "0.01011"
